$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns store plain-text values (e.g. "274.77", "1.001")
# rather than numbers. Force text format first so Excel's COM layer doesn't
# silently coerce numeric-looking strings (or thousands-separated "25.544.57"
# style values) into floating point numbers when we assign them below.
$ws.Range("D2:E51").NumberFormat = "@"

# Update cryptocurrency price/volume data (and re-sorted rows 35-51)
$ws.Range("D2").Value = '25.544.57'
$ws.Range("D3").Value = '1.809.09'
$ws.Range("E3").Value = '  -3.17%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '274.77'
$ws.Range("E5").Value = '  -8.46%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = '0.5012'
$ws.Range("E7").Value = '  -5.80%  '
$ws.Range("D8").Value = '0.3408'
$ws.Range("E8").Value = '  -8.58%  '
$ws.Range("D9").Value = '44.08'
$ws.Range("E9").Value = '  -2.81%  '
$ws.Range("D10").Value = '0.06633'
$ws.Range("E10").Value = '  -7.34%  '
$ws.Range("D11").Value = '19.48'
$ws.Range("E11").Value = '  -9.29%  '
$ws.Range("D12").Value = '0.7967'
$ws.Range("E12").Value = '  -10.09%  '
$ws.Range("D13").Value = '0.07841'
$ws.Range("E13").Value = '  -4.03%  '
$ws.Range("D14").Value = '1.810.39'
$ws.Range("E14").Value = '  -3.04%  '
$ws.Range("D15").Value = '5.018'
$ws.Range("E15").Value = '  -5.13%  '
$ws.Range("D16").Value = '86.63'
$ws.Range("E16").Value = '  -6.13%  '
$ws.Range("D17").Value = '1.0000'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '13.96'
$ws.Range("E18").Value = '  -5.80%  '
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("D20").Value = '0.000007937'
$ws.Range("E20").Value = '  -6.38%  '
$ws.Range("D21").Value = '25.605.13'
$ws.Range("E21").Value = '  -4.38%  '
$ws.Range("D22").Value = '4.704'
$ws.Range("E22").Value = '  -5.32%  '
$ws.Range("D23").Value = '9.858'
$ws.Range("E23").Value = '  -7.09%  '
$ws.Range("D24").Value = '6.088'
$ws.Range("E24").Value = '  -4.24%  '
$ws.Range("D25").Value = '2.246'
$ws.Range("E25").Value = '  -1.79%  '
$ws.Range("D26").Value = '142.41'
$ws.Range("E26").Value = '  -2.18%  '
$ws.Range("D27").Value = '1.659'
$ws.Range("E27").Value = '  -4.24%  '
$ws.Range("D28").Value = '17.02'
$ws.Range("E28").Value = '  -5.47%  '
$ws.Range("D29").Value = '108.19'
$ws.Range("E29").Value = '  -4.77%  '
$ws.Range("D30").Value = '4.254'
$ws.Range("E30").Value = '  -9.22%  '
$ws.Range("E31").Value = '  -9.47%  '
$ws.Range("D32").Value = '0.08689'
$ws.Range("E32").Value = '  -4.59%  '
$ws.Range("D33").Value = '0.04769'
$ws.Range("E33").Value = '  -4.78%  '
$ws.Range("D34").Value = '1.124'
$ws.Range("E34").Value = '  -3.85%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.7112'
$ws.Range("E35").Value = '  -11.48%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.833'
$ws.Range("E36").Value = '  -3.71%  '
$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").Value = '1.000'
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '3.111'
$ws.Range("E38").Value = '  -1.97%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '2.307'
$ws.Range("E39").Value = '  -13.27%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.01826'
$ws.Range("E40").Value = '  -5.85%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.5007'
$ws.Range("E41").Value = '  -17.90%  '
$ws.Range("D42").Value = '115.81'
$ws.Range("E42").Value = '  +0.94%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.9301'
$ws.Range("E43").Value = '  -12.40%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '6.143'
$ws.Range("E44").Value = '  -5.21%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '7.755'
$ws.Range("E46").Value = '  -11.06%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.1344'
$ws.Range("E47").Value = '  -9.85%  '
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = '0.4348'
$ws.Range("E48").Value = '  -16.62%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '35.98'
$ws.Range("E49").Value = '  -3.50%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '9.155'
$ws.Range("E50").Value = '  -7.44%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05896'
$ws.Range("E51").Value = '  -2.68%  '
